# Timesheet update: add/adjust daily hours booked against several tasks
# for the period covering columns L (day 9), O (day 12) and Q/R (days 14/15).
# All of the summary rows (Regular Hours / Over Time / Total Hrs / AI totals)
# are formula-driven (SUM / IF), so updating the task-level entry cells below
# is enough for the whole sheet to recompute consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - Unit Test Plan Preparation
$ws.Range("L14").Value = 1
$ws.Range("O14").Value = 1

# Row 19 - LLD Rework
$ws.Range("O19").Value = 1.5

# Row 24 - Code Peer Review
$ws.Range("Q24").Value = 1

# Row 25 - Peer Testing
$ws.Range("Q25").Value = 1

# Row 26 - Test Result review
$ws.Range("Q26").Value = 1

# Row 27 - Rework
$ws.Range("Q27").Value = 2

# Row 28 - Code Integration
$ws.Range("Q28").Value = 2

# Row 29 - Integration Testing
$ws.Range("Q29").Value = 1

# Row 30 - Rework
$ws.Range("Q30").Value = 1

# Row 31 - Test Result review
$ws.Range("Q31").Value = 1

# Row 32 - Presentation Preparation
$ws.Range("R32").Value = 3

# Row 33 - Deployment
$ws.Range("R33").Value = 1

# Leave the cursor where the author left it after the edits.
$ws.Range("S15").Select()
